$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = 1520
$ws.Range("E3").Value = 3012
$ws.Range("E5").Value = 114
$ws.Range("C8").Value = 597
$ws.Range("E9").Value = 113
$ws.Range("C12").Value = 409
$ws.Range("D12").Value = 719
$ws.Range("E12").Value = 1249
$ws.Range("C14").Value = 2632
$ws.Range("D14").Value = 5223
$ws.Range("E14").Value = 10380
$ws.Range("C15").Value = 199
$ws.Range("D15").Value = 370
$ws.Range("E15").Value = 722
$ws.Range("D18").Value = 64
$ws.Range("E18").Value = 126
$ws.Range("C20").Value = 652
$ws.Range("D20").Value = 1223
$ws.Range("E20").Value = 2300
$ws.Range("C24").Value = 830
$ws.Range("D24").Value = 1653
$ws.Range("E24").Value = 3326
$ws.Range("C27").Value = 160
$ws.Range("D27").Value = 319
$ws.Range("E27").Value = 648
$ws.Range("C29").Value = 10
$ws.Range("D29").Value = 29
$ws.Range("E29").Value = 72
$ws.Range("D30").Value = 1052
$ws.Range("E30").Value = 2119
$ws.Range("C31").Value = 743
$ws.Range("D31").Value = 1358
$ws.Range("E31").Value = 2668
$ws.Range("C37").Value = 2631
$ws.Range("D37").Value = 4629
$ws.Range("E37").Value = 8412
$ws.Range("C38").Value = 6
$ws.Range("C41").Value = 99
$ws.Range("D41").Value = 209
$ws.Range("E41").Value = 400
$ws.Range("C45").Value = 9
$ws.Range("D45").Value = 18
$ws.Range("E45").Value = 31
$ws.Range("C48").Value = 675
$ws.Range("D48").Value = 1309
$ws.Range("E48").Value = 2394
$ws.Range("D50").Value = 29
$ws.Range("E50").Value = 63
$ws.Range("C51").Value = 261
$ws.Range("D51").Value = 519
$ws.Range("E51").Value = 996
$ws.Range("D59").Value = 397
$ws.Range("E59").Value = 739
$ws.Range("C60").Value = 277
$ws.Range("D60").Value = 577
$ws.Range("E60").Value = 1151
$ws.Range("D63").Value = 261
$ws.Range("E66").Value = 98
$ws.Range("C67").Value = 326
$ws.Range("D67").Value = 666
$ws.Range("E67").Value = 1352
$ws.Range("C68").Value = 162
$ws.Range("D68").Value = 312
$ws.Range("E68").Value = 596
$ws.Range("C71").Value = 298
$ws.Range("D71").Value = 609
$ws.Range("E71").Value = 1266
$ws.Range("D72").Value = 459
$ws.Range("E72").Value = 969
$ws.Range("C73").Value = 169
$ws.Range("D73").Value = 340
$ws.Range("E73").Value = 699
$ws.Range("C79").Value = 72
$ws.Range("D79").Value = 156
$ws.Range("E79").Value = 314
$ws.Range("C80").Value = 88
$ws.Range("D80").Value = 184
$ws.Range("E80").Value = 389
$ws.Range("C81").Value = 218
$ws.Range("D81").Value = 432
$ws.Range("E81").Value = 853
$ws.Range("C82").Value = 2736
$ws.Range("D82").Value = 5176
$ws.Range("E82").Value = 9786
$ws.Range("C83").Value = 98
$ws.Range("D83").Value = 201
$ws.Range("E83").Value = 410
$ws.Range("C86").Value = 183
$ws.Range("D86").Value = 374
$ws.Range("E86").Value = 772
$ws.Range("D87").Value = 491
$ws.Range("E87").Value = 1024
$ws.Range("C90").Value = 427
$ws.Range("E90").Value = 1790
$ws.Range("D105").Value = 257
$ws.Range("E109").Value = 1536
